# Updated: pi 28. 05. 2021
# Apply new/revised AgTests (F) and AgPosit (G) values for rows 393-448,
# and append a new data row (449) for date 2021-05-27 (Excel serial 44343).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (F value, G value)
$updates = @{
    393 = @(308016, 1240)
    394 = @(165775, 623)
    395 = @(752222, 1950)
    396 = @(166408, 549)
    398 = @(298723, 1466)
    399 = @(200323, 967)
    400 = @(148957, 763)
    401 = @(272297, 932)
    402 = @(721834, 1393)
    403 = @(353822, 734)
    404 = @(224126, 914)
    405 = @(174024, 693)
    406 = @(170944, 680)
    407 = @(158075, 673)
    408 = @(304552, 835)
    409 = @(708107, 1006)
    410 = @(364399, 635)
    411 = @(225424, 828)
    412 = @(176157, 646)
    413 = @(149578, 658)
    414 = @(148839, 563)
    415 = @(307787, 694)
    416 = @(671568, 931)
    422 = @(298331, 645)
    423 = @(439321, 637)
    424 = @(266143, 499)
    425 = @(138000, 544)
    430 = @(175510, 271)
    435 = @(82389, 266)
    436 = @(145046, 352)
    437 = @(166850, 274)
    438 = @(121450, 251)
    440 = @(73394, 224)
    441 = @(68184, 202)
    442 = @(69967, 171)
    443 = @(106059, 208)
    444 = @(103104, 190)
    445 = @(84295, 189)
    446 = @(86215, 261)
    447 = @(67010, 192)
    448 = @(60904, 133)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
}

# Append the new row 449 with the latest day's data
$ws.Cells.Item(449, 1).Value = 44343
$ws.Cells.Item(449, 1).NumberFormat = $ws.Cells.Item(448, 1).NumberFormat
$ws.Cells.Item(449, 2).Value = 389440
$ws.Cells.Item(449, 3).Value = 4182
$ws.Cells.Item(449, 4).Value = 96
$ws.Cells.Item(449, 5).Value = 12333
$ws.Cells.Item(449, 6).Value = 43954
$ws.Cells.Item(449, 7).Value = 118
